# Auto commit at 2025-09-20  8:08:09.05
# Append two new daily rows (2025-09-19 / serial 45919) for both stations,
# following the same layout as the previous day's rows (37 -> 38/39).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (number formats/styles) of the last existing row down to
# the two new rows, same way Excel "fill down" would carry formatting.
$ws.Range("A37:F37").Copy() | Out-Null
$ws.Range("A38:F39").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Row 38: 四方坪站 (station "4")
$ws.Range("A38").Value = 45919
$ws.Range("B38").Value = "四方坪站"
$ws.Range("C38").Value = 7561.8
$ws.Range("D38").Value = 6053.53
$ws.Range("E38").Value = 2621.61
$ws.Range("F38").Value = 346

# Row 39: 高岭站 (station "5")
$ws.Range("A39").Value = 45919
$ws.Range("B39").Value = "高岭站"
$ws.Range("C39").Value = 4617.05
$ws.Range("D39").Value = 3656.83
$ws.Range("E39").Value = 1166.45
$ws.Range("F39").Value = 166

# Match the saved selection state recorded in the workbook after the edit.
$ws.Range("I38").Select() | Out-Null
